# Update the date heading and every "a OP b = c" answer in the table to the
# new values from the commit. Each call is a unique, case-sensitive, whole
# Find.Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,
#              MatchSoundsLike, MatchAllWordForms, Forward, Wrap,
#              Format, ReplaceWith, Replace)
# ReplaceWith text so a plain Find/ReplaceAll (wdReplaceAll = 2) is safe.
$d = $word.ActiveDocument
$d.Content.Find.Execute("2023-07-26 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-07-27 Thursday", 2) | Out-Null
$d.Content.Find.Execute("3+63=66", $true, $false, $false, $false, $false, $true, 1, $false, "57-11=46", 2) | Out-Null
$d.Content.Find.Execute("25+24=49", $true, $false, $false, $false, $false, $true, 1, $false, "97-34=63", 2) | Out-Null
$d.Content.Find.Execute("43-25=18", $true, $false, $false, $false, $false, $true, 1, $false, "15+10=25", 2) | Out-Null
$d.Content.Find.Execute("89-59=30", $true, $false, $false, $false, $false, $true, 1, $false, "38+2=40", 2) | Out-Null
$d.Content.Find.Execute("95-40=55", $true, $false, $false, $false, $false, $true, 1, $false, "33-12=21", 2) | Out-Null
$d.Content.Find.Execute("76-76=0", $true, $false, $false, $false, $false, $true, 1, $false, "40-19=21", 2) | Out-Null
$d.Content.Find.Execute("88-42=46", $true, $false, $false, $false, $false, $true, 1, $false, "15+37=52", 2) | Out-Null
$d.Content.Find.Execute("30-16=14", $true, $false, $false, $false, $false, $true, 1, $false, "40+44=84", 2) | Out-Null
$d.Content.Find.Execute("69+28=97", $true, $false, $false, $false, $false, $true, 1, $false, "1+74=75", 2) | Out-Null
$d.Content.Find.Execute("6+51=57", $true, $false, $false, $false, $false, $true, 1, $false, "69-34=35", 2) | Out-Null
$d.Content.Find.Execute("39-2=37", $true, $false, $false, $false, $false, $true, 1, $false, "74-29=45", 2) | Out-Null
$d.Content.Find.Execute("73+26=99", $true, $false, $false, $false, $false, $true, 1, $false, "21+22=43", 2) | Out-Null
$d.Content.Find.Execute("67-18=49", $true, $false, $false, $false, $false, $true, 1, $false, "24+29=53", 2) | Out-Null
$d.Content.Find.Execute("54-24=30", $true, $false, $false, $false, $false, $true, 1, $false, "35+18=53", 2) | Out-Null
$d.Content.Find.Execute("66-18=48", $true, $false, $false, $false, $false, $true, 1, $false, "37-3=34", 2) | Out-Null
$d.Content.Find.Execute("56+11=67", $true, $false, $false, $false, $false, $true, 1, $false, "75-21=54", 2) | Out-Null
$d.Content.Find.Execute("78-72=6", $true, $false, $false, $false, $false, $true, 1, $false, "88-48=40", 2) | Out-Null
$d.Content.Find.Execute("98-30=68", $true, $false, $false, $false, $false, $true, 1, $false, "99-34=65", 2) | Out-Null
$d.Content.Find.Execute("41+42=83", $true, $false, $false, $false, $false, $true, 1, $false, "17+68=85", 2) | Out-Null
$d.Content.Find.Execute("84-69=15", $true, $false, $false, $false, $false, $true, 1, $false, "91+8=99", 2) | Out-Null
$d.Content.Find.Execute("31-23=8", $true, $false, $false, $false, $false, $true, 1, $false, "16+51=67", 2) | Out-Null
$d.Content.Find.Execute("28+34=62", $true, $false, $false, $false, $false, $true, 1, $false, "63-47=16", 2) | Out-Null
$d.Content.Find.Execute("73-36=37", $true, $false, $false, $false, $false, $true, 1, $false, "67-52=15", 2) | Out-Null
$d.Content.Find.Execute("88-49=39", $true, $false, $false, $false, $false, $true, 1, $false, "29+61=90", 2) | Out-Null
$d.Content.Find.Execute("35+30=65", $true, $false, $false, $false, $false, $true, 1, $false, "57-4=53", 2) | Out-Null
$d.Content.Find.Execute("87+1=88", $true, $false, $false, $false, $false, $true, 1, $false, "65+31=96", 2) | Out-Null
$d.Content.Find.Execute("77-49=28", $true, $false, $false, $false, $false, $true, 1, $false, "16-7=9", 2) | Out-Null
$d.Content.Find.Execute("15+25=40", $true, $false, $false, $false, $false, $true, 1, $false, "64+4=68", 2) | Out-Null
$d.Content.Find.Execute("84-55=29", $true, $false, $false, $false, $false, $true, 1, $false, "35+15=50", 2) | Out-Null
$d.Content.Find.Execute("64-16=48", $true, $false, $false, $false, $false, $true, 1, $false, "80+9=89", 2) | Out-Null
$d.Content.Find.Execute("9+0=9", $true, $false, $false, $false, $false, $true, 1, $false, "35+45=80", 2) | Out-Null
$d.Content.Find.Execute("43-27=16", $true, $false, $false, $false, $false, $true, 1, $false, "51+36=87", 2) | Out-Null
$d.Content.Find.Execute("49-37=12", $true, $false, $false, $false, $false, $true, 1, $false, "78-53=25", 2) | Out-Null
$d.Content.Find.Execute("49-21=28", $true, $false, $false, $false, $false, $true, 1, $false, "73-49=24", 2) | Out-Null
$d.Content.Find.Execute("72-71=1", $true, $false, $false, $false, $false, $true, 1, $false, "87-10=77", 2) | Out-Null
$d.Content.Find.Execute("23+12=35", $true, $false, $false, $false, $false, $true, 1, $false, "51+21=72", 2) | Out-Null
$d.Content.Find.Execute("98-69=29", $true, $false, $false, $false, $false, $true, 1, $false, "51+11=62", 2) | Out-Null
$d.Content.Find.Execute("90-49=41", $true, $false, $false, $false, $false, $true, 1, $false, "26+42=68", 2) | Out-Null
$d.Content.Find.Execute("39-3=36", $true, $false, $false, $false, $false, $true, 1, $false, "90-55=35", 2) | Out-Null
$d.Content.Find.Execute("88-21=67", $true, $false, $false, $false, $false, $true, 1, $false, "81-53=28", 2) | Out-Null
$d.Content.Find.Execute("0+83=83", $true, $false, $false, $false, $false, $true, 1, $false, "21+60=81", 2) | Out-Null
$d.Content.Find.Execute("34+19=53", $true, $false, $false, $false, $false, $true, 1, $false, "34+11=45", 2) | Out-Null
$d.Content.Find.Execute("73-47=26", $true, $false, $false, $false, $false, $true, 1, $false, "94-3=91", 2) | Out-Null
$d.Content.Find.Execute("96-24=72", $true, $false, $false, $false, $false, $true, 1, $false, "18+13=31", 2) | Out-Null
$d.Content.Find.Execute("12+20=32", $true, $false, $false, $false, $false, $true, 1, $false, "67-41=26", 2) | Out-Null
$d.Content.Find.Execute("96-36=60", $true, $false, $false, $false, $false, $true, 1, $false, "36-27=9", 2) | Out-Null
$d.Content.Find.Execute("34+7=41", $true, $false, $false, $false, $false, $true, 1, $false, "84-59=25", 2) | Out-Null
$d.Content.Find.Execute("39+49=88", $true, $false, $false, $false, $false, $true, 1, $false, "16+74=90", 2) | Out-Null
$d.Content.Find.Execute("75-6=69", $true, $false, $false, $false, $false, $true, 1, $false, "24+42=66", 2) | Out-Null
$d.Content.Find.Execute("72+14=86", $true, $false, $false, $false, $false, $true, 1, $false, "44-35=9", 2) | Out-Null
$d.Content.Find.Execute("38+9=47", $true, $false, $false, $false, $false, $true, 1, $false, "21-0=21", 2) | Out-Null
$d.Content.Find.Execute("44-11=33", $true, $false, $false, $false, $false, $true, 1, $false, "6+42=48", 2) | Out-Null
$d.Content.Find.Execute("13-5=8", $true, $false, $false, $false, $false, $true, 1, $false, "68+5=73", 2) | Out-Null
$d.Content.Find.Execute("13+37=50", $true, $false, $false, $false, $false, $true, 1, $false, "15+0=15", 2) | Out-Null
$d.Content.Find.Execute("57+42=99", $true, $false, $false, $false, $false, $true, 1, $false, "57-27=30", 2) | Out-Null
$d.Content.Find.Execute("71-46=25", $true, $false, $false, $false, $false, $true, 1, $false, "96-70=26", 2) | Out-Null
$d.Content.Find.Execute("54+2=56", $true, $false, $false, $false, $false, $true, 1, $false, "90-52=38", 2) | Out-Null
$d.Content.Find.Execute("22+45=67", $true, $false, $false, $false, $false, $true, 1, $false, "76+5=81", 2) | Out-Null
$d.Content.Find.Execute("76-69=7", $true, $false, $false, $false, $false, $true, 1, $false, "15+62=77", 2) | Out-Null
$d.Content.Find.Execute("54+15=69", $true, $false, $false, $false, $false, $true, 1, $false, "38+26=64", 2) | Out-Null
$d.Content.Find.Execute("63+34=97", $true, $false, $false, $false, $false, $true, 1, $false, "66+5=71", 2) | Out-Null
$d.Content.Find.Execute("49+28=77", $true, $false, $false, $false, $false, $true, 1, $false, "78-16=62", 2) | Out-Null
$d.Content.Find.Execute("73-8=65", $true, $false, $false, $false, $false, $true, 1, $false, "15+64=79", 2) | Out-Null
$d.Content.Find.Execute("3+32=35", $true, $false, $false, $false, $false, $true, 1, $false, "75-7=68", 2) | Out-Null
$d.Content.Find.Execute("58-53=5", $true, $false, $false, $false, $false, $true, 1, $false, "75+14=89", 2) | Out-Null
$d.Content.Find.Execute("0+96=96", $true, $false, $false, $false, $false, $true, 1, $false, "20+13=33", 2) | Out-Null
$d.Content.Find.Execute("76-8=68", $true, $false, $false, $false, $false, $true, 1, $false, "67+27=94", 2) | Out-Null
$d.Content.Find.Execute("32-28=4", $true, $false, $false, $false, $false, $true, 1, $false, "40-2=38", 2) | Out-Null
$d.Content.Find.Execute("94-31=63", $true, $false, $false, $false, $false, $true, 1, $false, "88+0=88", 2) | Out-Null
$d.Content.Find.Execute("79-23=56", $true, $false, $false, $false, $false, $true, 1, $false, "23+68=91", 2) | Out-Null
$d.Content.Find.Execute("55-16=39", $true, $false, $false, $false, $false, $true, 1, $false, "83-70=13", 2) | Out-Null
$d.Content.Find.Execute("46+11=57", $true, $false, $false, $false, $false, $true, 1, $false, "11+61=72", 2) | Out-Null
$d.Content.Find.Execute("68+21=89", $true, $false, $false, $false, $false, $true, 1, $false, "42+23=65", 2) | Out-Null
$d.Content.Find.Execute("89-29=60", $true, $false, $false, $false, $false, $true, 1, $false, "65-62=3", 2) | Out-Null
$d.Content.Find.Execute("15+42=57", $true, $false, $false, $false, $false, $true, 1, $false, "60-26=34", 2) | Out-Null
$d.Content.Find.Execute("15+20=35", $true, $false, $false, $false, $false, $true, 1, $false, "10+82=92", 2) | Out-Null
$d.Content.Find.Execute("32+43=75", $true, $false, $false, $false, $false, $true, 1, $false, "32+41=73", 2) | Out-Null
$d.Content.Find.Execute("11+8=19", $true, $false, $false, $false, $false, $true, 1, $false, "68+2=70", 2) | Out-Null
$d.Content.Find.Execute("17+55=72", $true, $false, $false, $false, $false, $true, 1, $false, "76+14=90", 2) | Out-Null
$d.Content.Find.Execute("34+29=63", $true, $false, $false, $false, $false, $true, 1, $false, "24+48=72", 2) | Out-Null
$d.Content.Find.Execute("93-33=60", $true, $false, $false, $false, $false, $true, 1, $false, "75-19=56", 2) | Out-Null
$d.Content.Find.Execute("80-54=26", $true, $false, $false, $false, $false, $true, 1, $false, "62-11=51", 2) | Out-Null
$d.Content.Find.Execute("95-54=41", $true, $false, $false, $false, $false, $true, 1, $false, "70-30=40", 2) | Out-Null
$d.Content.Find.Execute("96-19=77", $true, $false, $false, $false, $false, $true, 1, $false, "26+18=44", 2) | Out-Null
$d.Content.Find.Execute("77-17=60", $true, $false, $false, $false, $false, $true, 1, $false, "23+37=60", 2) | Out-Null
$d.Content.Find.Execute("75+12=87", $true, $false, $false, $false, $false, $true, 1, $false, "20+19=39", 2) | Out-Null
$d.Content.Find.Execute("23-10=13", $true, $false, $false, $false, $false, $true, 1, $false, "77-67=10", 2) | Out-Null
$d.Content.Find.Execute("52+26=78", $true, $false, $false, $false, $false, $true, 1, $false, "96-93=3", 2) | Out-Null
$d.Content.Find.Execute("96-2=94", $true, $false, $false, $false, $false, $true, 1, $false, "27+34=61", 2) | Out-Null
$d.Content.Find.Execute("81-0=81", $true, $false, $false, $false, $false, $true, 1, $false, "93-38=55", 2) | Out-Null
$d.Content.Find.Execute("87-86=1", $true, $false, $false, $false, $false, $true, 1, $false, "64+17=81", 2) | Out-Null
$d.Content.Find.Execute("24+13=37", $true, $false, $false, $false, $false, $true, 1, $false, "78+5=83", 2) | Out-Null
$d.Content.Find.Execute("24+40=64", $true, $false, $false, $false, $false, $true, 1, $false, "73-35=38", 2) | Out-Null
$d.Content.Find.Execute("60-60=0", $true, $false, $false, $false, $false, $true, 1, $false, "90-23=67", 2) | Out-Null
$d.Content.Find.Execute("8+76=84", $true, $false, $false, $false, $false, $true, 1, $false, "32-16=16", 2) | Out-Null
$d.Content.Find.Execute("1+63=64", $true, $false, $false, $false, $false, $true, 1, $false, "81-77=4", 2) | Out-Null
$d.Content.Find.Execute("5+40=45", $true, $false, $false, $false, $false, $true, 1, $false, "82-52=30", 2) | Out-Null
$d.Content.Find.Execute("56-24=32", $true, $false, $false, $false, $false, $true, 1, $false, "92-23=69", 2) | Out-Null
$d.Content.Find.Execute("92-69=23", $true, $false, $false, $false, $false, $true, 1, $false, "15+83=98", 2) | Out-Null
$d.Content.Find.Execute("89-0=89", $true, $false, $false, $false, $false, $true, 1, $false, "11-6=5", 2) | Out-Null
